# Update the "想去人数" (interested count) values in column F for both the
# "展览" and "全部类型" worksheets, which carry duplicate data.
$wb = $excel.ActiveWorkbook

# Map of row -> new value, as per the commit's regenerated data export.
$updates = @{
    5  = 2619
    9  = 1361
    13 = 1176
    17 = 31
    19 = 72
    21 = 2461
    22 = 29
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
